# Generate Report for Archive
#
# 1) The localization status for zh-cn / de-de moved from "Ready for
#    handoff" to "In Translation" - update every cell that carried that
#    status string (the "Overview" summary sheet plus each language
#    sheet's own "Status" column).
# 2) The Status-related columns got narrower (the two status columns on
#    "Overview", and the "Status" column on each language sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the status columns. ColumnWidth is expressed in the standard
# Excel "characters" unit; 12.5 is what lands the stored column width on
# the target value used by the rest of the report.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
